# This script replicates the existing "Domaine Joseph Colin" product block
# (currently occupying rows 9-18) three more times into rows 19-48, so the
# wine list now contains four synchronized copies of the same catalogue
# entries (one original + three refreshed re-imports). The first line of
# each new block uses a slightly different label for the same "Aligoté"
# cuvee, matching two brand-new product names that need to be introduced
# into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used by the data rows: A Millesime, B Cuvee, C Domaine,
# D Appellation, G Pricetobuy, J Quantity, K Updated.
# (Columns E, F, H, I, L, M are intentionally left blank, matching the
# existing rows 9-18 this data is modeled on.)

$domaine     = "Domaine Joseph Colin"
$bourgogne   = "Bourgogne Générique"
$chassagne   = "Chassagne-Montrachet"
$puligny     = "Puligny-Montrachet"
$saintAubin  = "Saint-Aubin"
$updated     = 45684

# The ten-row template shared by every block (Millesime, Cuvee, Appellation, Pricetobuy).
# The Cuvee of the first row is supplied separately per block.
$template = @(
    @{ A = 0;    Cuvee = $null;                                 D = $bourgogne;  G = 15 },
    @{ A = 2021; Cuvee = "Chardonnay Les Hauts de la Combe";     D = $bourgogne;  G = 25 },
    @{ A = 2021; Cuvee = "Blanc";                                D = $chassagne;  G = 47 },
    @{ A = 2021; Cuvee = "En Cailleret";                         D = $chassagne;  G = 78 },
    @{ A = 2021; Cuvee = "Rouge Vieilles Vignes";                D = $chassagne;  G = 32 },
    @{ A = 2021; Cuvee = "La Garenne";                           D = $puligny;    G = 91 },
    @{ A = 2021; Cuvee = "Le Trezin 13";                         D = $puligny;    G = 85 },
    @{ A = 2021; Cuvee = "Compendium 135";                       D = $saintAubin; G = 32 },
    @{ A = 2021; Cuvee = "Clos du Meix 13";                      D = $saintAubin; G = 43 },
    @{ A = 2021; Cuvee = "La Chatenière 135";                    D = $saintAubin; G = 47 }
)

# First-row Cuvee label for each of the three new blocks (rows 19-28, 29-38, 39-48).
$firstRowLabels = @(
    "20Æaine Aligoté Les Jardins de la Cote",
    "Domaine Aligoté Les Jardins de la Cote",
    "Aligoté Les Jardins de la Cote"
)

$row = 19
foreach ($label in $firstRowLabels) {
    for ($i = 0; $i -lt $template.Count; $i++) {
        $data = $template[$i]

        if ($i -eq 0) {
            $cuvee = $label
        } else {
            $cuvee = $data.Cuvee
        }

        $ws.Cells.Item($row, 1).Value = $data.A
        $ws.Cells.Item($row, 2).Value = $cuvee
        $ws.Cells.Item($row, 3).Value = $domaine
        $ws.Cells.Item($row, 4).Value = $data.D
        $ws.Cells.Item($row, 7).Value = $data.G
        $ws.Cells.Item($row, 10).Value = 6
        $ws.Cells.Item($row, 11).Value = $updated

        $row = $row + 1
    }
}
